# Commit ECR and BEP(ECR is still unstable)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the ECR test case id text (B2) and wrap text to match other cells
$ws.Range("B2").Value = "ECR-10Apr2018-6"
$ws.Range("B2").WrapText = $true

# Move the active selection from F2 to B2
[void]$ws.Range("B2").Select()
